# Commit: Add files via upload
# - Rename sheet 'SFL' to 'SAF' and update the 'SFL' mentions in its facts to 'SAF'
# - Fix a typo in the QUE sheet's last fact ('than' -> 'that', trailing '?' -> '.')
# - Restore cursor/selection positions on the SAF and QUE sheets

$wb = $excel.ActiveWorkbook

# --- Rename SFL -> SAF, fix the airfield-code mentions in its facts ---
$saf = $wb.Worksheets.Item("SFL")
$saf.Name = "SAF"

$saf.Range("B2").Value = 'Over 1,700 flights landed in Vespuchia from the Saint Francois airfield, the country with the highest number of destined flights from SAF. There were almost double the amount of international flights than domestic flights out of this airfield.'
$saf.Range("B3").Value = 'The two highest destined cities from flights that departed from SAF were in Vespuchia and made up over 53% of departures.'
$saf.Range("B4").Value = 'Sundays were the busiest travel days for SAF and Tuesdays were the least busy travel days - if you’re looking to travel, we recommend traveling on Tuesdays for shorter lines and less crowded areas.'
$saf.Range("B6").Value = 'If you’re looking for shorter lines at security, we recommend visiting between 6 PM and 4 AM as there were no reported S2 scans during that time period at SAF.'
$saf.Range("B8").Value = 'The majority of passengers waited less than 16 minutes in line between S1 and S2 at the SAF airfield.'

# --- Fix the typo in the QUE sheet's last fact ---
$que = $wb.Worksheets.Item("QUE")
$que.Range("B6").Value = 'Did you know that there were over 100,000 passengers that landed in the Queenston airfield in 2028 and 2030 combined? And of those 100,000 only 379 of them occurred in 2030.'

# --- Restore the cell selections on each sheet ---
$saf.Activate()
$saf.Range("B21").Select()

$que.Activate()
$que.Range("B17").Select()

